$d = $word.ActiveDocument

$replacements = @(
    @("617×3=", "906×4="),
    @("833×8=", "345×9="),
    @("166×2=", "521×9="),
    @("407×3=", "867×2="),
    @("853×8=", "187×9="),
    @("863×3=", "907×8="),
    @("138×8=", "597×9="),
    @("499×3=", "293×7="),
    @("854×9=", "352×5="),
    @("207×6=", "447×3="),
    @("490×5=", "620×2="),
    @("703×4=", "388×4="),
    @("607×8=", "847×4="),
    @("510×4=", "867×5="),
    @("777×8=", "214×8="),
    @("522×5=", "354×9="),
    @("286×8=", "116×5="),
    @("946×7=", "352×6="),
    @("116×3=", "694×2="),
    @("155×9=", "983×2="),
    @("153×3=", "989×6="),
    @("330×8=", "450×5="),
    @("311×4=", "175×5="),
    @("770×4=", "431×3="),
    @("999×3=", "419×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
